# vr_gratings_params.xlsx - "Updates to run VR tuning experiments"
#
# - Insert a new column before the old "notes" column (F) and label it
#   "sortby"; the old "notes" column (and its data) shifts from F to G.
# - Tweak the isi value in B2 from "[0.00111, 0.04444]" to "[0.0111, 0.04444]".
# - Bump trial_duration (C2) and isi (E2) from 2/1 to 5/5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift column F ("notes" + "test") out to column G, freeing up F for the
# new "sortby" column.
$ws.Columns("F").Insert()

# New header cell for the inserted column.
$ws.Range("F1").Value = "sortby"

# Row 2 data tweaks.
$ws.Range("B2").Value = "[0.0111, 0.04444]"
$ws.Range("C2").Value = 5
$ws.Range("E2").Value = 5

# Give the new "sortby" column a bit more breathing room than a bare
# auto-fit would give it.
$ws.Columns("F").ColumnWidth = 6.75
